# Insert a new weekly record for "Macroferia Regional de Talca" - Apio,
# shifting the existing rows 304:345 down to 305:346, and populate the
# newly-opened row 304 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 304 (pushes 304:345 -> 305:346)
$ws.Rows.Item(304).Insert()

# Populate the new row 304 with the new record's values
$ws.Range("A304").Value = 5
$ws.Range("B304").Value = 'Macroferia Regional de Talca'
$ws.Range("C304").Value = 'Maule'
$ws.Range("D304").Value = 45142
$ws.Range("E304").Value = 7
$ws.Range("F304").Value = 100112017
$ws.Range("G304").Value = 'Apio'
$ws.Range("H304").Value = 'Americana (o)'
$ws.Range("I304").Value = 'Primera'
$ws.Range("J304").Value = 700
$ws.Range("K304").Value = 5000
$ws.Range("L304").Value = 5000
$ws.Range("M304").Value = 5000
$ws.Range("N304").Value = '$/docena de matas'
$ws.Range("O304").Value = 'Provincia del Elquí'
$ws.Range("P304").Value = 833
$ws.Range("Q304").Value = 6
$ws.Range("R304").Value = 'Hortaliza'
